# Add two new columns (I: "I0", J: "IF") to the worksheet, matching the
# header style already used by the other header cells (row 1) and filling
# in the numeric values for every data row (2-20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (font/border/alignment) of the existing header
# cell H1 onto the two new header cells so they look the same as the
# rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Re-apply the header text after the paste (PasteSpecial with formats
# only does not touch values, but keep this explicit just in case).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows ---------------------------------------------------------
$dataI = @(8, 8, 8, 8, 9, 8, 6, 7, 7, 5, 7, 6, 7, 8, 6, 5, 8, 5, 5)
$dataJ = @(8, 9, 9, 9, 9, 8, 7, 8, 8, 5, 7, 6, 7, 8, 7, 5, 8, 5, 5)

for ($i = 0; $i -lt $dataI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}

Write-Output "Added columns I and J (I0 / IF) with header + data for rows 1-20"
